$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from an existing header
# cell so the new header cells (AD1:AF1) pick up the same style index
# instead of minting a new one.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels for the season-record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
$ws.Range("AD2:AD47").Value = 88
$ws.Range("AE2:AE47").Value = 74
$ws.Range("AF2:AF47").Value = 0
